$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill Summary")

# Helper: write a numeric-looking value into a cell while keeping it stored
# as TEXT (matching the workbook's existing convention of quoting amount
# cells as strings, e.g. "20480.00"). Assigning such a string straight to
# .Value would be auto-coerced back into a real number by Excel, so instead
# we briefly place a text formula, then paste-special just the value over
# itself - this "bakes" the literal text result without touching any cell
# formatting/styles.
function Set-TextNumber([string]$addr, [string]$text) {
    $cell = $ws.Range($addr)
    $cell.Formula = "=""$text"""
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# "Qty executed upto date" (column C) - plain numbers
$ws.Range("C8").Value = 69
$ws.Range("C9").Value = 82
$ws.Range("C10").Value = 91
$ws.Range("C11").Value = 33
$ws.Range("C12").Value = 97
$ws.Range("C13").Value = 56
$ws.Range("C14").Value = 29
$ws.Range("C15").Value = 85
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 68

# "Upto date Amount" (column G), stays text e.g. "20992.00"
Set-TextNumber "G9"  "20992.00"
Set-TextNumber "G10" "42952.00"
Set-TextNumber "G11" "21846.00"
Set-TextNumber "G13" "7616.00"
Set-TextNumber "G14" "667.00"

# Grand-total rows
Set-TextNumber "G19" "94073.00"
Set-TextNumber "H19" "94073.00"
Set-TextNumber "G21" "94073.00"
Set-TextNumber "H21" "94073.00"
